# Update workbook to reflect data through 2022-08-12 (adds one more day of data)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab and the sheet title to reflect the new "through" date
$ws.Name = "Through 2022-08-12"

# Update the August row label text
$ws.Range("A9").Value = "August (through 08-12)"

# Update August row (row 9) values for years 2015-2022 (columns B-I)
$ws.Range("B9").Value = 12
$ws.Range("C9").Value = 25
$ws.Range("D9").Value = 26
$ws.Range("E9").Value = 21
$ws.Range("F9").Value = 16
$ws.Range("G9").Value = 77
$ws.Range("H9").Value = 75
$ws.Range("I9").Value = 66

# Update Total row (row 10) values for years 2015-2022 (columns B-I)
$ws.Range("B10").Value = 174
$ws.Range("C10").Value = 327
$ws.Range("D10").Value = 491
$ws.Range("E10").Value = 446
$ws.Range("F10").Value = 320
$ws.Range("G10").Value = 698
$ws.Range("H10").Value = 985
$ws.Range("I10").Value = 1036
